$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency data for rows 2-51 (Coin, Link, Price, Volume(1h))
# Source data scraped Wed Jun  7 07:49:59 UTC 2023
$data = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "26.859.26", "  +4.20%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.878.64", "  +3.49%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.000", "  -0.10%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "277.63", "  +0.15%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.000", "  -0.08%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.5282", "  +4.11%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3416", "  -3.07%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06951", "  +4.34%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "20.06", "  +0.09%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.8039", "  -2.85%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07717", "  -1.89%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.882.87", "  +3.72%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "90.33", "  +3.26%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.179", "  +2.04%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "14.60", "  +3.36%  "),
    @("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "0.9994", "  -0.15%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000008071", "  +0.33%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9999", "  -0.06%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "26.905.82", "  +4.14%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.122.67", "  +4.22%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.744", "  +0.21%  "),
    @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "10.04", "  +0.28%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.199", "  +1.71%  "),
    @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.396", "  +9.24%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "146.36", "  +2.97%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.663", "  -0.62%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "17.31", "  +1.21%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "113.31", "  +3.59%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.354", "  +0.31%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.315", "  +1.86%  "),
    @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.08901", "  +1.21%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04899", "  +0.39%  "),
    @("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.175", "  +3.21%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7275", "  -0.20%  "),
    @("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.878", "  +0.28%  "),
    @("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "3.289", "  +5.10%  "),
    @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.351", "  -1.02%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01843", "  -0.37%  "),
    @("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.5130", "  -0.87%  "),
    @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.9560", "  -0.90%  "),
    @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "116.34", "  +5.16%  "),
    @("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.175", "  -0.30%  "),
    @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "8.099", "  +1.04%  "),
    @("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9996", "  -0.11%  "),
    @("Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.4461", "  -2.40%  "),
    @("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1340", "  -1.72%  "),
    @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.333", "  +0.48%  "),
    @("Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "36.21", "  -1.01%  "),
    @("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05946", "  +1.80%  "),
    @("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.491", "  -0.50%  ")
)

# Force column D (Price) to Text format so numeric-looking strings (e.g. "1.000",
# "14.60", "0.000008071") keep their exact printed form instead of being
# auto-coerced into numbers by the Value setter.
$ws.Range("D2:D51").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rowData[0]
    $ws.Cells.Item($row, 3).Value = $rowData[1]
    $ws.Cells.Item($row, 4).Value = $rowData[2]
    $ws.Cells.Item($row, 5).Value = $rowData[3]
}
